$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 431, pushing existing rows 431:492 down to 432:493.
$ws.Range("A431").EntireRow.Insert()

# Populate the newly inserted row 431 with the new record's data.
$ws.Range("A431").Value = 3
$ws.Range("B431").Value = "Femacal de La Calera"
$ws.Range("C431").Value = "Coquimbo"
$ws.Range("D431").Value = 45131
$ws.Range("D431").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E431").Value = 5
$ws.Range("F431").Value = 100112001
$ws.Range("G431").Value = "Berenjena"
$ws.Range("H431").Value = "Sin especificar"
$ws.Range("I431").Value = "Primera"
$ws.Range("J431").Value = 130
$ws.Range("K431").Value = 7000
$ws.Range("L431").Value = 7500
$ws.Range("M431").Value = 7250
$ws.Range("N431").Value = "`$/caja 60 unidades"
$ws.Range("O431").Value = "Región de Arica y Parinacota"
$ws.Range("P431").Value = 121
$ws.Range("Q431").Value = 60
$ws.Range("R431").Value = "Hortaliza"
